# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", before "2022-Q2".
# 2. Populate "2022-Q3" with the new fund-holdings table (007139 / 011556).
# 3. Insert a new row into "总计" for the "2022-Q3" quarter, pushing the
#    existing 2022-Q2 / 2021-Q2 / 2020-Q4 rows down by one.
#
# xlPasteFormats / xlPasteValues constants used below (Copy+PasteSpecial)
# so that numeric-looking text (fund codes, percentages, ...) round-trips
# as plain text without leaving stray number-format styles or formulas
# behind, and so that new cells inherit the same style as their
# neighbours.
$xlPasteValues = -4163
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet right after "总计" (tab position 2)
# ------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q2_2022 = $wb.Worksheets.Item("2022-Q2")

$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Mirror the header / index-column styling from the "2022-Q2" sheet.
$q2_2022.Range("B1:H1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats) | Out-Null
$q2_2022.Range("A2:A3").Copy() | Out-Null
$q3.Range("A2:A3").PasteSpecial($xlPasteFormats) | Out-Null

# Header row.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2 - 007139
$q3.Range("A2").Value = 0
$q3.Range("C2").Value = "富国民裕进取沪港深成长精选混合A"
$q3.Range("H2").Value = 10

# Row 3 - 011556
$q3.Range("A3").Value = 1
$q3.Range("C3").Value = "富国民裕进取沪港深成长精选混合C"
$q3.Range("H3").Value = 10

# Numeric-looking values must stay plain text (t="inlineStr"/shared
# string), matching fund codes like "007139" that have significant
# leading zeros and decimal strings like "10.56" that must not be
# reformatted as doubles. Enter them as formulas first, then convert
# the whole block to static values in one shot so no leftover <f>
# survives and no custom number-format style gets attached.
$q3.Range("B2").Formula = '="007139"'
$q3.Range("D2").Formula = '="10.56"'
$q3.Range("E2").Formula = '="88.09"'
$q3.Range("F2").Formula = '="4.71"'
$q3.Range("G2").Formula = '="0.4974"'

$q3.Range("B3").Formula = '="011556"'
$q3.Range("D3").Formula = '="2.46"'
$q3.Range("E3").Formula = '="88.09"'
$q3.Range("F3").Formula = '="4.71"'
$q3.Range("G3").Formula = '="0.1159"'

$q3text = $q3.Range("B2:G3")
$q3text.Copy() | Out-Null
$q3text.PasteSpecial($xlPasteValues) | Out-Null

# ------------------------------------------------------------------
# 2) Insert the "2022-Q3" row into "总计", shifting the rest down.
# ------------------------------------------------------------------
# Extend formatting down to the new last row (row 5) before filling it,
# by cloning the style of the (soon to be) previous last data row.
$zongji.Range("A4").Copy() | Out-Null
$zongji.Range("A5").PasteSpecial($xlPasteFormats) | Out-Null

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2020-Q4"
$zongji.Range("C5").Value = 6
$zongji.Range("D5").Value = 1.22

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q2"
$zongji.Range("C4").Value = 2
$zongji.Range("D4").Value = 0.44

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q2"
$zongji.Range("C3").Value = 2
$zongji.Range("D3").Value = 0.24

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.61
